$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$profs = @(
    'Osvaldo Boyle PhD',
    'Darwin Nicolas',
    'Mr. Howard Willms II',
    'Carissa Adams',
    'Prof. Eladio Franecki',
    'Dr. Yvonne Treutel PhD',
    'Laisha Schultz',
    'Maymie Haley',
    'Destinee Feest',
    'Vernie Sporer',
)

$periods = @(
    'Wednesday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Thursday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
    'Sunday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Sunday-12:30PM:2:00PM,Sunday-2:00PM:3:30PM,Sunday-3:30PM:5:00PM,Monday-9:00AM:10:30AM,Monday-10:30AM:12:00PM,Monday-12:30PM:2:00PM,Monday-2:00PM:3:30PM,Monday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM,Tuesday-10:30AM:12:00PM,Tuesday-12:30PM:2:00PM,Tuesday-2:00PM:3:30PM,Tuesday-3:30PM:5:00PM,Wednesday-9:00AM:10:30AM,Wednesday-10:30AM:12:00PM,Wednesday-12:30PM:2:00PM,Wednesday-2:00PM:3:30PM,Wednesday-3:30PM:5:00PM,Thursday-9:00AM:10:30AM,Thursday-10:30AM:12:00PM,Thursday-12:30PM:2:00PM,Thursday-2:00PM:3:30PM,Thursday-3:30PM:5:00PM,Saturday-9:00AM:10:30AM,Saturday-10:30AM:12:00PM,Saturday-12:30PM:2:00PM,Saturday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
    'Sunday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Sunday-12:30PM:2:00PM,Sunday-2:00PM:3:30PM,Sunday-3:30PM:5:00PM,Monday-9:00AM:10:30AM,Monday-10:30AM:12:00PM,Monday-12:30PM:2:00PM,Monday-2:00PM:3:30PM,Monday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM,Tuesday-10:30AM:12:00PM,Tuesday-12:30PM:2:00PM,Tuesday-2:00PM:3:30PM,Tuesday-3:30PM:5:00PM,Wednesday-9:00AM:10:30AM,Wednesday-10:30AM:12:00PM,Wednesday-12:30PM:2:00PM,Wednesday-2:00PM:3:30PM,Wednesday-3:30PM:5:00PM,Thursday-9:00AM:10:30AM,Thursday-10:30AM:12:00PM,Thursday-12:30PM:2:00PM,Thursday-2:00PM:3:30PM,Thursday-3:30PM:5:00PM,Saturday-9:00AM:10:30AM,Saturday-10:30AM:12:00PM,Saturday-12:30PM:2:00PM,Saturday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
    'Sunday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Sunday-12:30PM:2:00PM,Sunday-2:00PM:3:30PM,Sunday-3:30PM:5:00PM,Monday-9:00AM:10:30AM,Monday-10:30AM:12:00PM,Monday-12:30PM:2:00PM,Monday-2:00PM:3:30PM,Monday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM,Tuesday-10:30AM:12:00PM,Tuesday-12:30PM:2:00PM,Tuesday-2:00PM:3:30PM,Tuesday-3:30PM:5:00PM,Wednesday-9:00AM:10:30AM,Wednesday-10:30AM:12:00PM,Wednesday-12:30PM:2:00PM,Wednesday-2:00PM:3:30PM,Wednesday-3:30PM:5:00PM,Thursday-9:00AM:10:30AM,Thursday-10:30AM:12:00PM,Thursday-12:30PM:2:00PM,Thursday-2:00PM:3:30PM,Thursday-3:30PM:5:00PM,Saturday-9:00AM:10:30AM,Saturday-10:30AM:12:00PM,Saturday-12:30PM:2:00PM,Saturday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
    'Sunday-10:30AM:12:00PM,Sunday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM',
    'Sunday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Sunday-12:30PM:2:00PM,Sunday-2:00PM:3:30PM,Sunday-3:30PM:5:00PM,Monday-9:00AM:10:30AM,Monday-10:30AM:12:00PM,Monday-12:30PM:2:00PM,Monday-2:00PM:3:30PM,Monday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM,Tuesday-10:30AM:12:00PM,Tuesday-12:30PM:2:00PM,Tuesday-2:00PM:3:30PM,Tuesday-3:30PM:5:00PM,Wednesday-9:00AM:10:30AM,Wednesday-10:30AM:12:00PM,Wednesday-12:30PM:2:00PM,Wednesday-2:00PM:3:30PM,Wednesday-3:30PM:5:00PM,Thursday-9:00AM:10:30AM,Thursday-10:30AM:12:00PM,Thursday-12:30PM:2:00PM,Thursday-2:00PM:3:30PM,Thursday-3:30PM:5:00PM,Saturday-9:00AM:10:30AM,Saturday-10:30AM:12:00PM,Saturday-12:30PM:2:00PM,Saturday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
    'Sunday-10:30AM:12:00PM',
    'Sunday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Sunday-12:30PM:2:00PM,Sunday-2:00PM:3:30PM,Sunday-3:30PM:5:00PM,Monday-9:00AM:10:30AM,Monday-10:30AM:12:00PM,Monday-12:30PM:2:00PM,Monday-2:00PM:3:30PM,Monday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM,Tuesday-10:30AM:12:00PM,Tuesday-12:30PM:2:00PM,Tuesday-2:00PM:3:30PM,Tuesday-3:30PM:5:00PM,Wednesday-9:00AM:10:30AM,Wednesday-10:30AM:12:00PM,Wednesday-12:30PM:2:00PM,Wednesday-2:00PM:3:30PM,Wednesday-3:30PM:5:00PM,Thursday-9:00AM:10:30AM,Thursday-10:30AM:12:00PM,Thursday-12:30PM:2:00PM,Thursday-2:00PM:3:30PM,Thursday-3:30PM:5:00PM,Saturday-9:00AM:10:30AM,Saturday-10:30AM:12:00PM,Saturday-12:30PM:2:00PM,Saturday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
    'Sunday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Sunday-12:30PM:2:00PM,Sunday-2:00PM:3:30PM,Sunday-3:30PM:5:00PM,Monday-9:00AM:10:30AM,Monday-10:30AM:12:00PM,Monday-12:30PM:2:00PM,Monday-2:00PM:3:30PM,Monday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM,Tuesday-10:30AM:12:00PM,Tuesday-12:30PM:2:00PM,Tuesday-2:00PM:3:30PM,Tuesday-3:30PM:5:00PM,Wednesday-9:00AM:10:30AM,Wednesday-10:30AM:12:00PM,Wednesday-12:30PM:2:00PM,Wednesday-2:00PM:3:30PM,Wednesday-3:30PM:5:00PM,Thursday-9:00AM:10:30AM,Thursday-10:30AM:12:00PM,Thursday-12:30PM:2:00PM,Thursday-2:00PM:3:30PM,Thursday-3:30PM:5:00PM,Saturday-9:00AM:10:30AM,Saturday-10:30AM:12:00PM,Saturday-12:30PM:2:00PM,Saturday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
    'Sunday-9:00AM:10:30AM,Sunday-10:30AM:12:00PM,Sunday-12:30PM:2:00PM,Sunday-2:00PM:3:30PM,Sunday-3:30PM:5:00PM,Monday-9:00AM:10:30AM,Monday-10:30AM:12:00PM,Monday-12:30PM:2:00PM,Monday-2:00PM:3:30PM,Monday-3:30PM:5:00PM,Tuesday-9:00AM:10:30AM,Tuesday-10:30AM:12:00PM,Tuesday-12:30PM:2:00PM,Tuesday-2:00PM:3:30PM,Tuesday-3:30PM:5:00PM,Wednesday-9:00AM:10:30AM,Wednesday-10:30AM:12:00PM,Wednesday-12:30PM:2:00PM,Wednesday-2:00PM:3:30PM,Wednesday-3:30PM:5:00PM,Thursday-9:00AM:10:30AM,Thursday-10:30AM:12:00PM,Thursday-12:30PM:2:00PM,Thursday-2:00PM:3:30PM,Thursday-3:30PM:5:00PM,Saturday-9:00AM:10:30AM,Saturday-10:30AM:12:00PM,Saturday-12:30PM:2:00PM,Saturday-2:00PM:3:30PM,Saturday-3:30PM:5:00PM',
)

$startRow = 25
for ($i = 0; $i -lt $profs.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $profs[$i]
    $ws.Cells.Item($r, 2).Value = $periods[$i]
}
